$d = $word.ActiveDocument

$pairs = @(
    @("13×89=1157", "21×82=1722"),
    @("46×83=3818", "22×21=462"),
    @("40×46=1840", "82×99=8118"),
    @("17×12=204", "65×31=2015"),
    @("75×30=2250", "41×71=2911"),
    @("28×68=1904", "11×44=484"),
    @("56×61=3416", "58×48=2784"),
    @("13×90=1170", "52×42=2184"),
    @("64×72=4608", "47×80=3760"),
    @("54×54=2916", "47×15=705"),
    @("42×78=3276", "12×62=744"),
    @("19×31=589", "96×75=7200"),
    @("64×87=5568", "50×81=4050"),
    @("94×11=1034", "76×57=4332"),
    @("67×11=737", "98×94=9212"),
    @("92×92=8464", "72×61=4392"),
    @("52×90=4680", "71×94=6674"),
    @("47×81=3807", "46×19=874"),
    @("51×92=4692", "47×29=1363"),
    @("94×13=1222", "50×74=3700"),
    @("47×99=4653", "24×28=672"),
    @("15×81=1215", "49×56=2744"),
    @("18×66=1188", "21×66=1386"),
    @("99×87=8613", "66×40=2640"),
    @("72×52=3744", "48×34=1632")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
